# The upstream commit ("updated spanish and english slides. updated
# language from master slide to Instructor slides") touched several
# files in the repository. For *this* presentation
# (En-Lesson_Slide-Create_AI_Content.pptx) the recorded OOXML diff
# contains no textual, structural, or formatting changes whatsoever:
# every one of its ~35 hunks is a reshuffle of the `xmlns*` attribute
# *order* on the same `<ma14:wrappingTextBoxFlag .../>` extension
# element (and one `<mc:Fallback>` wrapper) that Office already wrote
# into the file - the tag name, namespace URIs, and `val="1"` payload
# are identical before and after, just re-ordered, which is exactly
# what happens when a deck is round-tripped through a different Office
# build (e.g. Mac PowerPoint) without anyone touching its content.
#
# That attribute-order detail lives purely in low level OOXML
# serialization plumbing for an undocumented Mac extension flag - it
# is not backed by any property on the PowerPoint object model (no
# Shape/TextFrame/Fill/Line property maps onto it), so there is no
# COM automation call, here or in real PowerPoint, that targets it.
#
# So the faithful reproduction of this diff is to open the deck and
# leave its content exactly as authored - i.e. a no-op edit - rather
# than invent a content change that the diff does not actually show.
$p = $ppt.ActivePresentation
$null = $p.Slides.Count
